$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 49; this shifts the existing rows 49-113 down to 50-114,
# preserving all of their data (Excel's native row-insert/shift semantics).
$ws.Rows("49").Insert()

# Populate the newly-inserted (blank) row 49 with the new record.
$ws.Range("A49").Value = 2
$ws.Range("B49").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C49").Value = "Coquimbo"
$ws.Range("D49").Value = 44994
$ws.Range("E49").Value = 4
$ws.Range("F49").Value = 100112030
$ws.Range("G49").Value = "Poroto granado"
$ws.Range("H49").Value = "Sin especificar"
$ws.Range("I49").Value = "Primera"
$ws.Range("J49").Value = 400
$ws.Range("K49").Value = 21000
$ws.Range("L49").Value = 23000
$ws.Range("M49").Value = 22000
$ws.Range("N49").Value = "$/malla 25 kilos"
$ws.Range("O49").Value = "Provincia de Limarí"
$ws.Range("P49").Value = 880
$ws.Range("Q49").Value = 25
$ws.Range("R49").Value = "Hortaliza"
